$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("AI9").Value = 11
$ws.Range("AU9").Value = 7.5
$ws.Range("G9").Value = 2.3
$ws.Range("I9").Value = 2.8
$ws.Range("J9").Value = 2.88
$ws.Range("L9").Value = 3.4

# Row 17
$ws.Range("M17").Value = 1.06
$ws.Range("N17").Value = 10
$ws.Range("O17").Value = 1.33
$ws.Range("P17").Value = 3.25

# Row 28
$ws.Range("AA28").Value = 21
$ws.Range("AC28").Value = 7.3
$ws.Range("AE28").Value = 13.5
$ws.Range("AK28").Value = 20
$ws.Range("AL28").Value = 28
$ws.Range("AP28").Value = 23
$ws.Range("AT28").Value = 2.67
$ws.Range("AU28").Value = 7.1
$ws.Range("AX28").Value = 13.5
$ws.Range("G28").Value = 2.55
$ws.Range("H28").Value = 3.25
$ws.Range("I28").Value = 2.5
$ws.Range("N28").Value = 7.3
$ws.Range("Q28").Value = 1.9
$ws.Range("T28").Value = 2.67
$ws.Range("W28").Value = 8.5
$ws.Range("Z28").Value = 28

# Row 32
$ws.Range("AA32").Value = 20
$ws.Range("AD32").Value = 6.3
$ws.Range("AE32").Value = 13
$ws.Range("AG32").Value = 9.25
$ws.Range("AH32").Value = 14
$ws.Range("AI32").Value = 9.75
$ws.Range("AJ32").Value = 30
$ws.Range("AN32").Value = 4.5
$ws.Range("AO32").Value = 13.5
$ws.Range("AU32").Value = 6.9
$ws.Range("AV32").Value = 60
$ws.Range("AW32").Value = 4.6
$ws.Range("G32").Value = 2.47
$ws.Range("H32").Value = 3.2
$ws.Range("I32").Value = 2.6
$ws.Range("J32").Value = 3.15
$ws.Range("K32").Value = 2.1
$ws.Range("L32").Value = 3.2
$ws.Range("Q32").Value = 1.87
$ws.Range("R32").Value = 1.87
$ws.Range("W32").Value = 8.5
$ws.Range("X32").Value = 13
$ws.Range("Z32").Value = 27

# Row 34
$ws.Range("AC34").Value = 7.1
$ws.Range("AI34").Value = 9.5
$ws.Range("AJ34").Value = 26
$ws.Range("AO34").Value = 14
$ws.Range("AR34").Value = 100
$ws.Range("AT34").Value = 2.75
$ws.Range("G34").Value = 2.57
$ws.Range("H34").Value = 3.3
$ws.Range("J34").Value = 3.2
$ws.Range("K34").Value = 2.12
$ws.Range("N34").Value = 7.1
$ws.Range("Q34").Value = 1.98
$ws.Range("S34").Value = 1.4
$ws.Range("T34").Value = 2.75
$ws.Range("X34").Value = 13

# Row 36
$ws.Range("AA36").Value = 25
$ws.Range("AI36").Value = 8.75
$ws.Range("AJ36").Value = 24
$ws.Range("AK36").Value = 18
$ws.Range("AL36").Value = 25
$ws.Range("AN36").Value = 4.9
$ws.Range("AO36").Value = 16.5
$ws.Range("AP36").Value = 23
$ws.Range("AQ36").Value = 80
$ws.Range("AR36").Value = 110
$ws.Range("AS36").Value = 300
$ws.Range("AT36").Value = 2.52
$ws.Range("AU36").Value = 6.6
$ws.Range("AW36").Value = 4.25
$ws.Range("AX36").Value = 12
$ws.Range("AY36").Value = 19
$ws.Range("BA36").Value = 75
$ws.Range("G36").Value = 2.95
$ws.Range("I36").Value = 2.3
$ws.Range("J36").Value = 3.55
$ws.Range("K36").Value = 2.05
$ws.Range("L36").Value = 2.87
$ws.Range("O36").Value = 1.27
$ws.Range("P36").Value = 3.1
$ws.Range("S36").Value = 1.39
$ws.Range("T36").Value = 2.55
$ws.Range("X36").Value = 16
$ws.Range("Y36").Value = 10.5
